$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace "conb_car" / "E_car" labels in column B with "PC_petrol" / "PC_electric"
# (every other row, starting at row 2: conb_car on even rows, E_car on odd rows)
for ($r = 2; $r -le 17; $r++) {
    $val = $ws.Cells.Item($r, 2).Value()
    if ($val -eq "conb_car") {
        $ws.Cells.Item($r, 2).Value = "PC_petrol"
    } elseif ($val -eq "E_car") {
        $ws.Cells.Item($r, 2).Value = "PC_electric"
    }
}

# Update the active selection to B3
$ws.Range("B3").Select()
